$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings like "42.802.38" / "316.60" that must stay
# literal text. Assigning a plain decimal string (e.g. "316.26") via .Value
# would otherwise be auto-coerced to a Number (and drop a trailing zero, as
# in "105.80" -> 105.8). Force text format right before the write, then reset
# the cell style back to Normal so no stray style/format is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '42.868.23'
$ws.Range("E2").Value = '  +0.39%  '
Set-TextValue $ws.Range("D3") '2.530.55'
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("E4").Value = '  +0.02%  '
Set-TextValue $ws.Range("D5") '316.26'
$ws.Range("E5").Value = '  +1.18%  '
Set-TextValue $ws.Range("D6") '97.64'
$ws.Range("E6").Value = '  +2.49%  '
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -0.60%  '
Set-TextValue $ws.Range("D10") '35.84'
$ws.Range("E10").Value = '  -0.81%  '
Set-TextValue $ws.Range("D11") '0.0812'
$ws.Range("E11").Value = '  +0.32%  '
Set-TextValue $ws.Range("D12") '7.61'
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("E13").Value = '  -2.53%  '
Set-TextValue $ws.Range("D14") '2.912.83'
$ws.Range("E14").Value = '  +0.27%  '
Set-TextValue $ws.Range("D15") '2.547.23'
$ws.Range("E15").Value = '  +1.95%  '
Set-TextValue $ws.Range("D16") '15.16'
$ws.Range("E16").Value = '  -2.98%  '
$ws.Range("E17").Value = '  -1.43%  '
Set-TextValue $ws.Range("D18") '42.884.46'
$ws.Range("E18").Value = '  +0.37%  '
Set-TextValue $ws.Range("D19") '6.86'
$ws.Range("E19").Value = '  +4.91%  '
Set-TextValue $ws.Range("D20") '12.74'
$ws.Range("E20").Value = '  -2.86%  '
Set-TextValue $ws.Range("D21") '0.0₃0963'
$ws.Range("E21").Value = '  -0.50%  '
Set-TextValue $ws.Range("D22") '69.63'
$ws.Range("E22").Value = '  -2.30%  '
Set-TextValue $ws.Range("D23") '252.57'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("E25").Value = '  +0.10%  '
Set-TextValue $ws.Range("D26") '26.42'
$ws.Range("E26").Value = '  -1.97%  '
$ws.Range("E27").Value = '  -0.15%  '
$ws.Range("E28").Value = '  +2.77%  '
Set-TextValue $ws.Range("D29") '41.29'
$ws.Range("E29").Value = '  +5.08%  '
Set-TextValue $ws.Range("D30") '10.41'
$ws.Range("E30").Value = '  +3.67%  '
$ws.Range("E31").Value = '  +0.36%  '
Set-TextValue $ws.Range("D32") '158.51'
$ws.Range("E32").Value = '  +0.93%  '
$ws.Range("E33").Value = '  +3.58%  '
$ws.Range("E34").Value = '  +4.10%  '
Set-TextValue $ws.Range("D35") '3.33'
$ws.Range("E35").Value = '  +0.17%  '
Set-TextValue $ws.Range("D36") '18.92'
$ws.Range("E36").Value = '  -4.71%  '
Set-TextValue $ws.Range("D37") '0.0788'
$ws.Range("E37").Value = '  +0.21%  '
$ws.Range("E38").Value = '  -0.22%  '
Set-TextValue $ws.Range("D39") '2.48'
$ws.Range("E39").Value = '  +18.07%  '
$ws.Range("E40").Value = '  -0.84%  '
Set-TextValue $ws.Range("D41") '21.76'
$ws.Range("E41").Value = '  -10.89%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E43").Value = '  +1.35%  '
$ws.Range("E44").Value = '  +0.19%  '
Set-TextValue $ws.Range("D45") '3.31'
$ws.Range("E45").Value = '  -1.80%  '
Set-TextValue $ws.Range("D46") '2.018.32'
$ws.Range("E46").Value = '  -2.71%  '
Set-TextValue $ws.Range("D47") '9.07'
$ws.Range("E47").Value = '  +2.86%  '
Set-TextValue $ws.Range("D48") '84.17'
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue $ws.Range("D49") '75.88'
$ws.Range("E49").Value = '  +2.03%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D50") '105.80'
$ws.Range("E50").Value = '  +4.24%  '
Set-TextValue $ws.Range("D51") '2.772.31'
$ws.Range("E51").Value = '  +0.50%  '
